$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-driven decay values in column B (rows 2-6) with
# plain literal numbers (formulas removed, values kept).
$ws.Range("B2").Value = [double]"1.5779141105123101E-2"
$ws.Range("B3").Value = [double]"6.1024417699829603E-5"
$ws.Range("B4").Value = [double]"1.2068579208866E-6"
$ws.Range("B5").Value = [double]"3.0258892967883001E-8"
$ws.Range("B6").Value = [double]"6.07355146616586E-9"

# Clear out the rest of the old cycle-number / css-criterion data
# (rows 7-19 previously held A-column cycle numbers and B-column shared
# formulas); leave the cell formatting (style) in place.
$ws.Range("A7:A19").ClearContents()
$ws.Range("B7:B19").ClearContents()

# Remove the now-unused trailing blank rows (18-30) so the sheet shrinks
# back down to A1:B17.
$ws.Range("A18:B30").EntireRow.Delete()

# Move the active selection.
[void]$ws.Range("D8").Select()
